$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 197.5433703333333
$ws.Cells.Item(2, 8).Value = 592.6301109999999
$ws.Cells.Item(2, 9).Value = 0.3388703761585983
$ws.Cells.Item(2, 10).Value = 0.3388703761585982
$ws.Cells.Item(2, 15).Value = 0.7140239834365498
$ws.Cells.Item(2, 16).Value = 0.7140239834365498
$ws.Cells.Item(2, 17).Value = 395.5636762104414
$ws.Cells.Item(2, 18).Value = 3560.073085893973
$ws.Cells.Item(2, 19).Value = 0.2419615758534044
$ws.Cells.Item(2, 20).Value = 0.2419615758534043

# Row 3
$ws.Cells.Item(3, 7).Value = 197.5433703333333
$ws.Cells.Item(3, 8).Value = 592.6301109999999
$ws.Cells.Item(3, 9).Value = 0.3388703761585983
$ws.Cells.Item(3, 10).Value = 0.3388703761585982
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.8019933333333333
$ws.Cells.Item(3, 14).Value = 2.40598
$ws.Cells.Item(3, 15).Value = 0.2859760165634502
$ws.Cells.Item(3, 16).Value = 0.2859760165634502
$ws.Cells.Item(3, 17).Value = 158.4284660515311
$ws.Cells.Item(3, 18).Value = 1425.85619446378
$ws.Cells.Item(3, 19).Value = 0.09690880030519389
$ws.Cells.Item(3, 20).Value = 0.09690880030519387

# Row 4
$ws.Cells.Item(4, 9).Value = 0.1369374790620155
$ws.Cells.Item(4, 10).Value = 0.1369374790620154
$ws.Cells.Item(4, 15).Value = 0.7140239834365498
$ws.Cells.Item(4, 16).Value = 0.7140239834365498
$ws.Cells.Item(4, 19).Value = 0.09777664428161943
$ws.Cells.Item(4, 20).Value = 0.0977766442816194

# Row 5
$ws.Cells.Item(5, 9).Value = 0.1369374790620155
$ws.Cells.Item(5, 10).Value = 0.1369374790620154
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8019933333333333
$ws.Cells.Item(5, 14).Value = 2.40598
$ws.Cells.Item(5, 15).Value = 0.2859760165634502
$ws.Cells.Item(5, 16).Value = 0.2859760165634502
$ws.Cells.Item(5, 17).Value = 64.02092445698223
$ws.Cells.Item(5, 18).Value = 576.18832011284
$ws.Cells.Item(5, 19).Value = 0.03916083478039605
$ws.Cells.Item(5, 20).Value = 0.03916083478039604

# Row 6
$ws.Cells.Item(6, 7).Value = 148.824417
$ws.Cells.Item(6, 8).Value = 446.473251
$ws.Cells.Item(6, 9).Value = 0.2552967790580629
$ws.Cells.Item(6, 10).Value = 0.2552967790580629
$ws.Cells.Item(6, 15).Value = 0.7140239834365498
$ws.Cells.Item(6, 16).Value = 0.7140239834365498
$ws.Cells.Item(6, 17).Value = 298.008145750777
$ws.Cells.Item(6, 18).Value = 2682.073311756993
$ws.Cells.Item(6, 19).Value = 0.1822880231415588
$ws.Cells.Item(6, 20).Value = 0.1822880231415588

# Row 7
$ws.Cells.Item(7, 7).Value = 148.824417
$ws.Cells.Item(7, 8).Value = 446.473251
$ws.Cells.Item(7, 9).Value = 0.2552967790580629
$ws.Cells.Item(7, 10).Value = 0.2552967790580629
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.8019933333333333
$ws.Cells.Item(7, 14).Value = 2.40598
$ws.Cells.Item(7, 15).Value = 0.2859760165634502
$ws.Cells.Item(7, 16).Value = 0.2859760165634502
$ws.Cells.Item(7, 17).Value = 119.35619027122
$ws.Cells.Item(7, 18).Value = 1074.20571244098
$ws.Cells.Item(7, 19).Value = 0.07300875591650406
$ws.Cells.Item(7, 20).Value = 0.07300875591650406

# Row 8
$ws.Cells.Item(8, 7).Value = 35.426853
$ws.Cells.Item(8, 8).Value = 106.280559
$ws.Cells.Item(8, 9).Value = 0.06077202683121193
$ws.Cells.Item(8, 10).Value = 0.06077202683121192
$ws.Cells.Item(8, 15).Value = 0.7140239834365498
$ws.Cells.Item(8, 16).Value = 0.7140239834365498
$ws.Cells.Item(8, 17).Value = 70.939238232093
$ws.Cells.Item(8, 18).Value = 638.453144088837
$ws.Cells.Item(8, 19).Value = 0.04339268467953483
$ws.Cells.Item(8, 20).Value = 0.04339268467953482

# Row 9
$ws.Cells.Item(9, 7).Value = 35.426853
$ws.Cells.Item(9, 8).Value = 106.280559
$ws.Cells.Item(9, 9).Value = 0.06077202683121193
$ws.Cells.Item(9, 10).Value = 0.06077202683121192
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.8019933333333333
$ws.Cells.Item(9, 14).Value = 2.40598
$ws.Cells.Item(9, 15).Value = 0.2859760165634502
$ws.Cells.Item(9, 16).Value = 0.2859760165634502
$ws.Cells.Item(9, 17).Value = 28.41209992698
$ws.Cells.Item(9, 18).Value = 255.70889934282
$ws.Cells.Item(9, 19).Value = 0.0173793421516771
$ws.Cells.Item(9, 20).Value = 0.0173793421516771

# Row 10
$ws.Cells.Item(10, 7).Value = 121.3248153333333
$ws.Cells.Item(10, 8).Value = 363.974446
$ws.Cells.Item(10, 9).Value = 0.2081233388901116
$ws.Cells.Item(10, 10).Value = 0.2081233388901115
$ws.Cells.Item(10, 15).Value = 0.7140239834365498
$ws.Cells.Item(10, 16).Value = 0.7140239834365498
$ws.Cells.Item(10, 17).Value = 242.9425492124864
$ws.Cells.Item(10, 18).Value = 2186.482942912378
$ws.Cells.Item(10, 19).Value = 0.1486050554804325
$ws.Cells.Item(10, 20).Value = 0.1486050554804324

# Row 11
$ws.Cells.Item(11, 7).Value = 121.3248153333333
$ws.Cells.Item(11, 8).Value = 363.974446
$ws.Cells.Item(11, 9).Value = 0.2081233388901116
$ws.Cells.Item(11, 10).Value = 0.2081233388901115
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.8019933333333333
$ws.Cells.Item(11, 14).Value = 2.40598
$ws.Cells.Item(11, 15).Value = 0.2859760165634502
$ws.Cells.Item(11, 16).Value = 0.2859760165634502
$ws.Cells.Item(11, 17).Value = 97.3016930652311
$ws.Cells.Item(11, 18).Value = 875.71523758708
$ws.Cells.Item(11, 19).Value = 0.05951828340967909
$ws.Cells.Item(11, 20).Value = 0.05951828340967909
